$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

# Clear the obsolete pcap-name labels in column D (rows 5-8); chart only
# references columns E and F so these text labels are no longer needed.
$ws.Range("D5:D8").ClearContents()

# Update the active selection to match the author's final cursor position.
$ws.Range("K22").Select()
